{"js": "// Replace each three-digit-by-one-digit multiplication prompt in the\n// worksheet table with the newly generated problem, one-for-one.\nconst replacements = [\n  [\"671\u00d78=\", \"840\u00d77=\"],\n  [\"688\u00d72=\", \"993\u00d76=\"],\n  [\"395\u00d74=\", \"110\u00d75=\"],\n  [\"180\u00d77=\", \"488\u00d77=\"],\n  [\"191\u00d77=\", \"815\u00d76=\"],\n  [\"103\u00d76=\", \"215\u00d79=\"],\n  [\"194\u00d74=\", \"576\u00d73=\"],\n  [\"956\u00d79=\", \"555\u00d75=\"],\n  [\"330\u00d72=\", \"406\u00d74=\"],\n  [\"703\u00d78=\", \"254\u00d76=\"],\n  [\"981\u00d75=\", \"402\u00d74=\"],\n  [\"441\u00d72=\", \"846\u00d78=\"],\n  [\"725\u00d75=\", \"687\u00d76=\"],\n  [\"723\u00d79=\", \"671\u00d77=\"],\n  [\"615\u00d73=\", \"621\u00d75=\"],\n  [\"925\u00d73=\", \"838\u00d78=\"],\n  [\"155\u00d77=\", \"936\u00d77=\"],\n  [\"201\u00d77=\", \"682\u00d74=\"],\n  [\"444\u00d75=\", \"987\u00d74=\"],\n  [\"259\u00d79=\", \"757\u00d78=\"],\n  [\"414\u00d76=\", \"974\u00d78=\"],\n  [\"670\u00d76=\", \"438\u00d77=\"],\n  [\"607\u00d77=\", \"331\u00d77=\"],\n  [\"230\u00d74=\", \"850\u00d75=\"],\n  [\"568\u00d77=\", \"665\u00d79=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication prompt in the\n# worksheet table with the newly generated problem, one-for-one.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"671\u00d78=\", \"840\u00d77=\"),\n    @(\"688\u00d72=\", \"993\u00d76=\"),\n    @(\"395\u00d74=\", \"110\u00d75=\"),\n    @(\"180\u00d77=\", \"488\u00d77=\"),\n    @(\"191\u00d77=\", \"815\u00d76=\"),\n    @(\"103\u00d76=\", \"215\u00d79=\"),\n    @(\"194\u00d74=\", \"576\u00d73=\"),\n    @(\"956\u00d79=\", \"555\u00d75=\"),\n    @(\"330\u00d72=\", \"406\u00d74=\"),\n    @(\"703\u00d78=\", \"254\u00d76=\"),\n    @(\"981\u00d75=\", \"402\u00d74=\"),\n    @(\"441\u00d72=\", \"846\u00d78=\"),\n    @(\"725\u00d75=\", \"687\u00d76=\"),\n    @(\"723\u00d79=\", \"671\u00d77=\"),\n    @(\"615\u00d73=\", \"621\u00d75=\"),\n    @(\"925\u00d73=\", \"838\u00d78=\"),\n    @(\"155\u00d77=\", \"936\u00d77=\"),\n    @(\"201\u00d77=\", \"682\u00d74=\"),\n    @(\"444\u00d75=\", \"987\u00d74=\"),\n    @(\"259\u00d79=\", \"757\u00d78=\"),\n    @(\"414\u00d76=\", \"974\u00d78=\"),\n    @(\"670\u00d76=\", \"438\u00d77=\"),\n    @(\"607\u00d77=\", \"331\u00d77=\"),\n    @(\"230\u00d74=\", \"850\u00d75=\"),\n    @(\"568\u00d77=\", \"665\u00d79=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
